# Staging.ResultArea.xlsx - reproduce the tracked edit:
#   - header cells "ResultArea_ID" (A2) and "BusinessKey" (D2) are swapped
#   - (best effort / not exposed via Excel automation, see notes below)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the A2 / D2 header labels ("ResultArea_ID" <-> "BusinessKey") ----
$resultAreaId = $ws.Range("A2").Value()
$businessKey  = $ws.Range("D2").Value()

$ws.Range("A2").Value = $businessKey
$ws.Range("D2").Value = $resultAreaId

# --- Best-effort: restore the window size recorded in bookViews --------
# (xWindow/yWindow/windowWidth/windowHeight are not reachable through the
# documented Worksheet/Window COM surface in this runtime, but attempting
# the assignment is harmless if unsupported.)
try {
    $win = $wb.Windows.Item(1)
    $win.Width = 28800
    $win.Height = 12585
} catch {
}

# --- Best-effort: the worksheet's internal VBA CodeName ------------------
# (CodeName is a read-only pseudo-property on the automation object model -
# real Excel only lets the VBE rename it - so this is a no-op if rejected.)
try {
    $ws.CodeName = "Sheet47"
} catch {
}
